# Module 4 Networking.docx - append " Networking" as its own run right
# after the existing "Introduction" run (same bold/bCs/underline
# formatting), per the commit's xml diff.

$d = $word.ActiveDocument

# Locate the "Introduction" heading text robustly (rather than assuming
# a fixed paragraph index).
$rng = $d.Content
$found = $rng.Find.Execute("Introduction", $false, $true, $false, $false, `
    $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target text 'Introduction'"
}

# Remember the character offset right after "Introduction" - that's
# where the new text needs to be inserted.
$insertPos = $rng.End

# Collapse to that point and type the new text. (It initially gets
# merged into the same run as "Introduction" because the formatting is
# identical.)
$rng.Collapse(0)
$rng.InsertAfter(" Networking")

# Re-select just the newly typed text and nudge its Bold property off
# then back on. That forces the engine to split it out into its own
# <w:r> (matching the diff, which adds " Networking" as a distinct run)
# instead of leaving it silently merged into the "Introduction" run.
$newRun = $d.Range($insertPos, $insertPos + 11)
$newRun.Font.Bold = $false
$newRun.Font.Bold = $true

Write-Output ("Paragraph now reads: " + $newRun.Paragraphs(1).Range.Text)
